$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 12:39"

# --- Country re-ranking (shared-string / label swaps) ---
# Emiratos Arabes Unidos overtakes Guatemala and China
$ws.Cells.Item(44, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(45, 1).Value = "Guatemala"
$ws.Cells.Item(46, 1).Value = "China"

# Timor Oriental overtakes Santa Lucia
$ws.Cells.Item(204, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 1).Value = "Santa Lucia"

# Islas Malvinas overtakes Montserrat
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 1).Value = "Montserrat"

# --- Updated case numbers (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 18 - Banglades
$ws.Cells.Item(18, 2).Value = 350621
$ws.Cells.Item(18, 3).Value = 1705
$ws.Cells.Item(18, 4).Value = 258717
$ws.Cells.Item(18, 5).Value = 86925
$ws.Cells.Item(18, 7).Value = 40
$ws.Cells.Item(18, 8).Value = 4979

# Row 27 - Israel
$ws.Cells.Item(27, 2).Value = 188760
$ws.Cells.Item(27, 3).Value = 858
$ws.Cells.Item(27, 4).Value = 135997
$ws.Cells.Item(27, 5).Value = 51503
$ws.Cells.Item(27, 7).Value = 4
$ws.Cells.Item(27, 8).Value = 1260

# Row 33 - Rumania
$ws.Cells.Item(33, 2).Value = 113589
$ws.Cells.Item(33, 3).Value = 808
$ws.Cells.Item(33, 4).Value = 90649
$ws.Cells.Item(33, 5).Value = 18482
$ws.Cells.Item(33, 7).Value = 23
$ws.Cells.Item(33, 8).Value = 4458

# Row 44 - Emiratos Arabes Unidos (new leader of this group)
$ws.Cells.Item(44, 2).Value = 85595
$ws.Cells.Item(44, 3).Value = 679
$ws.Cells.Item(44, 4).Value = 75086
$ws.Cells.Item(44, 5).Value = 10104
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 405

# Row 45 - Guatemala
$ws.Cells.Item(45, 2).Value = 85444
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 74859
$ws.Cells.Item(45, 5).Value = 7466
$ws.Cells.Item(45, 8).Value = 3119

# Row 46 - China
$ws.Cells.Item(46, 2).Value = 85291
$ws.Cells.Item(46, 3).Value = 12
$ws.Cells.Item(46, 4).Value = 80484
$ws.Cells.Item(46, 5).Value = 173
$ws.Cells.Item(46, 8).Value = 4634

# Row 99 - Malasia
$ws.Cells.Item(99, 2).Value = 10276
$ws.Cells.Item(99, 3).Value = 57
$ws.Cells.Item(99, 4).Value = 9395
$ws.Cells.Item(99, 5).Value = 751

# Row 141
$ws.Cells.Item(141, 4).Value = 3100
$ws.Cells.Item(141, 5).Value = 174

# Row 145 - Sri Lanka
$ws.Cells.Item(145, 2).Value = 2776
$ws.Cells.Item(145, 3).Value = 45
$ws.Cells.Item(145, 4).Value = 2079
$ws.Cells.Item(145, 5).Value = 675
$ws.Cells.Item(145, 7).Value = 2
$ws.Cells.Item(145, 8).Value = 22

# Row 214 - Islas Malvinas (after relabel)
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

# Row 215 - Montserrat (after relabel)
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
